# Fruta / hortaliza, semanal
# The weekly refresh re-sorted/shuffled the daily price rows (2-31) for
# this sub-grouping; every row's variable fields (date, variety, quality,
# volume, prices, unit, origin, $/Kg, Kg-or-units) were replaced by the
# corresponding fields of a different source row. The "key" columns
# (A,B,C,E,F,G,R) are constant across all rows, so only columns D,H,I,J,K,
# L,M,N,O,P,Q need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of new-row -> old-row (1-indexed worksheet rows) that reproduces
# the post-edit ordering of the data block (rows 2..31).
$rowMap = [ordered]@{
    2  = 12
    3  = 5
    4  = 21
    5  = 7
    6  = 11
    7  = 28
    8  = 17
    9  = 30
    10 = 4
    11 = 16
    12 = 22
    13 = 23
    14 = 24
    15 = 8
    16 = 25
    17 = 13
    18 = 19
    19 = 20
    20 = 2
    21 = 9
    22 = 10
    23 = 27
    24 = 15
    25 = 31
    26 = 14
    27 = 29
    28 = 18
    29 = 26
    30 = 6
    31 = 3
}

# Columns whose values vary row-to-row and therefore need to be
# re-shuffled according to $rowMap.
$cols = 4,8,9,10,11,12,13,14,15,16,17

# 1. Snapshot the current (pre-edit) values of every cell we might touch,
#    BEFORE any writes happen, so that source rows used later as a
#    "from" row are unaffected by earlier writes.
$snapshot = @{}
foreach ($r in 2..31) {
    foreach ($c in $cols) {
        $snapshot["$r,$c"] = $ws.Cells.Item($r, $c).Value()
    }
}

# 2. Write the shuffled values back out.
foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($newRow, $c).Value = $snapshot["$oldRow,$c"]
    }
}
